$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46071.01041666666, 2463.053),
    @(3, 46071.02083333334, 2513.594),
    @(4, 46071.03125, 2512.251),
    @(5, 46071.04166666666, 2512.578),
    @(6, 46071.05208333334, 2509.834),
    @(7, 46071.0625, 2509.987),
    @(8, 46071.07291666666, 2513.228),
    @(9, 46071.08333333334, 2513.173),
    @(10, 46071.09375, 2503.831),
    @(11, 46071.10416666666, 2498.452),
    @(12, 46071.11458333334, 2493.971),
    @(13, 46071.125, 2489.608),
    @(14, 46071.13541666666, 2486.004),
    @(15, 46071.14583333334, 2479.585),
    @(16, 46071.15625, 2473.81),
    @(17, 46071.16666666666, 2468.688),
    @(18, 46071.17708333334, 2414.813),
    @(19, 46071.1875, 2396.077),
    @(20, 46071.19791666666, 2435.363),
    @(21, 46071.20833333334, 2425.5),
    @(22, 46071.21875, 2407.195),
    @(23, 46071.22916666666, 2396.278),
    @(24, 46071.23958333334, 2381.012),
    @(25, 46071.25, 2430.486),
    @(26, 46071.26041666666, 2346.001),
    @(27, 46071.27083333334, 2408.635),
    @(28, 46071.28125, 2397.472),
    @(29, 46071.29166666666, 2387.696),
    @(30, 46071.30208333334, 2368.854),
    @(31, 46071.3125, 2319.774),
    @(32, 46071.32291666666, 2282.062),
    @(33, 46071.33333333334, 2275.967),
    @(34, 46071.34375, 2266.263),
    @(35, 46071.35416666666, 2334.075),
    @(36, 46071.36458333334, 2330.932),
    @(37, 46071.375, 2287.097),
    @(38, 46071.38541666666, 2327.294),
    @(39, 46071.39583333334, 2327.533),
    @(40, 46071.40625, 2338.478),
    @(41, 46071.41666666666, 2271.155),
    @(42, 46071.42708333334, 2256.558),
    @(43, 46071.4375, 2268.089),
    @(44, 46071.44791666666, 2266.102),
    @(45, 46071.45833333334, 2263.024),
    @(46, 46071.46875, 2263.705),
    @(47, 46071.47916666666, 2261.263),
    @(48, 46071.48958333334, 2259.501),
    @(49, 46071.5, 2257.677),
    @(50, 46071.51041666666, 2252.241),
    @(51, 46071.52083333334, 2252.011),
    @(52, 46071.53125, 2249.037),
    @(53, 46071.54166666666, 2249.674),
    @(54, 46071.55208333334, 2245.542),
    @(55, 46071.5625, 2240.777),
    @(56, 46071.57291666666, 2236.367),
    @(57, 46071.58333333334, 2231.553),
    @(58, 46071.59375, 2234.336),
    @(59, 46071.60416666666, 2224.506),
    @(60, 46071.61458333334, 2214.764),
    @(61, 46071.625, 2206.038),
    @(62, 46071.63541666666, 2187.782),
    @(63, 46071.64583333334, 2176.793),
    @(64, 46071.65625, 2232.889),
    @(65, 46071.66666666666, 2220.198),
    @(66, 46071.67708333334, 2126.174),
    @(67, 46071.6875, 2106.884),
    @(68, 46071.69791666666, 2152.414),
    @(69, 46071.70833333334, 2131.175),
    @(70, 46071.71875, 2114.181),
    @(71, 46071.72916666666, 2083.805),
    @(72, 46071.73958333334, 2052.356),
    @(73, 46071.75, 2021.27),
    @(74, 46071.76041666666, 1978.173),
    @(75, 46071.77083333334, 1942.199),
    @(76, 46071.78125, 1905.566),
    @(77, 46071.79166666666, 1869.225),
    @(78, 46071.80208333334, 1823.235),
    @(79, 46071.8125, 1794.612),
    @(80, 46071.82291666666, 1765.596),
    @(81, 46071.83333333334, 1737.247),
    @(82, 46071.84375, 1682.928),
    @(83, 46071.85416666666, 1651.05),
    @(84, 46071.86458333334, 1618.691),
    @(85, 46071.875, 1586.118),
    @(86, 46071.88541666666, 1527.538),
    @(87, 46071.89583333334, 1493.436),
    @(88, 46071.90625, 1457.91),
    @(89, 46071.91666666666, 1423.211),
    @(90, 46071.92708333334, 1364.647),
    @(91, 46071.9375, 1324.051),
    @(92, 46071.94791666666, 1283.083),
    @(93, 46071.95833333334, 1243.139),
    @(94, 46071.96875, $null),
    @(95, 46071.97916666666, $null),
    @(96, 46071.98958333334, $null),
    @(97, 46072, $null)
)

foreach ($item in $data) {
    $row = $item[0]
    $aVal = $item[1]
    $bVal = $item[2]
    $ws.Cells.Item($row, 1).Value = $aVal
    if ($bVal -ne $null) {
        $ws.Cells.Item($row, 2).Value = $bVal
    }
}
